{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"2024-10-16 Wednesday\", \"2024-10-17 Thursday\"],\n  [\"84\u00d737=3108\", \"34\u00d734=1156\"],\n  [\"33\u00d797=3201\", \"94\u00d797=9118\"],\n  [\"17\u00d794=1598\", \"68\u00d772=4896\"],\n  [\"51\u00d796=4896\", \"82\u00d779=6478\"],\n  [\"61\u00d744=2684\", \"57\u00d745=2565\"],\n  [\"76\u00d737=2812\", \"11\u00d743=473\"],\n  [\"49\u00d716=784\", \"63\u00d782=5166\"],\n  [\"67\u00d723=1541\", \"15\u00d724=360\"],\n  [\"87\u00d748=4176\", \"77\u00d746=3542\"],\n  [\"57\u00d721=1197\", \"40\u00d732=1280\"],\n  [\"41\u00d744=1804\", \"22\u00d724=528\"],\n  [\"88\u00d721=1848\", \"35\u00d716=560\"],\n  [\"69\u00d764=4416\", \"78\u00d790=7020\"],\n  [\"41\u00d750=2050\", \"99\u00d760=5940\"],\n  [\"77\u00d785=6545\", \"36\u00d768=2448\"],\n  [\"30\u00d786=2580\", \"91\u00d792=8372\"],\n  [\"80\u00d773=5840\", \"37\u00d733=1221\"],\n  [\"18\u00d783=1494\", \"28\u00d795=2660\"],\n  [\"17\u00d762=1054\", \"13\u00d794=1222\"],\n  [\"55\u00d787=4785\", \"45\u00d750=2250\"],\n  [\"36\u00d746=1656\", \"20\u00d784=1680\"],\n  [\"84\u00d769=5796\", \"83\u00d775=6225\"],\n  [\"91\u00d727=2457\", \"66\u00d752=3432\"],\n  [\"86\u00d763=5418\", \"15\u00d791=1365\"],\n  [\"32\u00d768=2176\", \"42\u00d774=3108\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @{Old=\"2024-10-16 Wednesday\"; New=\"2024-10-17 Thursday\"},\n  @{Old=\"84\u00d737=3108\"; New=\"34\u00d734=1156\"},\n  @{Old=\"33\u00d797=3201\"; New=\"94\u00d797=9118\"},\n  @{Old=\"17\u00d794=1598\"; New=\"68\u00d772=4896\"},\n  @{Old=\"51\u00d796=4896\"; New=\"82\u00d779=6478\"},\n  @{Old=\"61\u00d744=2684\"; New=\"57\u00d745=2565\"},\n  @{Old=\"76\u00d737=2812\"; New=\"11\u00d743=473\"},\n  @{Old=\"49\u00d716=784\"; New=\"63\u00d782=5166\"},\n  @{Old=\"67\u00d723=1541\"; New=\"15\u00d724=360\"},\n  @{Old=\"87\u00d748=4176\"; New=\"77\u00d746=3542\"},\n  @{Old=\"57\u00d721=1197\"; New=\"40\u00d732=1280\"},\n  @{Old=\"41\u00d744=1804\"; New=\"22\u00d724=528\"},\n  @{Old=\"88\u00d721=1848\"; New=\"35\u00d716=560\"},\n  @{Old=\"69\u00d764=4416\"; New=\"78\u00d790=7020\"},\n  @{Old=\"41\u00d750=2050\"; New=\"99\u00d760=5940\"},\n  @{Old=\"77\u00d785=6545\"; New=\"36\u00d768=2448\"},\n  @{Old=\"30\u00d786=2580\"; New=\"91\u00d792=8372\"},\n  @{Old=\"80\u00d773=5840\"; New=\"37\u00d733=1221\"},\n  @{Old=\"18\u00d783=1494\"; New=\"28\u00d795=2660\"},\n  @{Old=\"17\u00d762=1054\"; New=\"13\u00d794=1222\"},\n  @{Old=\"55\u00d787=4785\"; New=\"45\u00d750=2250\"},\n  @{Old=\"36\u00d746=1656\"; New=\"20\u00d784=1680\"},\n  @{Old=\"84\u00d769=5796\"; New=\"83\u00d775=6225\"},\n  @{Old=\"91\u00d727=2457\"; New=\"66\u00d752=3432\"},\n  @{Old=\"86\u00d763=5418\"; New=\"15\u00d791=1365\"},\n  @{Old=\"32\u00d768=2176\"; New=\"42\u00d774=3108\"},\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair.Old\n  $find.Replacement.Text = $pair.New\n  $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n\nWrite-Output \"done\""}
